$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Objekt_ID" label to use a hyphen instead of an underscore
$ws.Range("A2").Value = "Objekt-ID"

# Wrap the title value in literal double quotes
$ws.Range("B3").Value = '"Zigeunerpaar"'

# Drop the stray formatted-but-empty row all the way at the bottom of the
# sheet (row 1048576) that was inflating the used range / dimension.
$ws.Rows.Item(1048576).Delete()

# Move / set the active selection to the last used cell (B32)
$ws.Range("B32").Select()
